# Edit script: updates the "VdG/C3/22" annex document.
#  1) Change the suscrito-con-fecha date from 31 to 26 of May 2022 in the
#     intro paragraph.
#  2) Remove the light-grey single-line cell borders from every cell of the
#     7-column student table (3 rows x 7 cols).
#  3) Update a couple of placeholder cell values and the start/end dates in
#     both data rows of that table.
#  4) Change the "En Puertollano a 31 de mayo 2022" signature date to 26.

$d = $word.ActiveDocument

# --- 1) Intro paragraph date -------------------------------------------------
$d.Content.Find.Execute(
    "suscrito con fecha  31 de mayo de 2022",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "suscrito con fecha  26 de mayo de 2022", 2)

# --- 2) Strip the single/888888 borders on every cell of the student table --
$t = $d.Tables.Item(2)
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    for ($c = 1; $c -le $t.Columns.Count; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Borders.Item(-1).LineStyle = 0
        $cell.Borders.Item(-2).LineStyle = 0
        $cell.Borders.Item(-3).LineStyle = 0
        $cell.Borders.Item(-4).LineStyle = 0
    }
}

# --- 3) Update the two data rows --------------------------------------------
# Row 2 (Díez Viñas Malena)
$t.Cell(2, 4).Range.Text = "szsc"
$t.Cell(2, 6).Range.Text = "2022-05-26"
$t.Cell(2, 7).Range.Text = "2022-05-27"

# Row 3 (Moreno Ramos Laura)
$t.Cell(3, 4).Range.Text = "zxcxc"
$t.Cell(3, 6).Range.Text = "2022-05-26"
$t.Cell(3, 7).Range.Text = "2022-05-27"

# --- 4) Signature line date --------------------------------------------------
$d.Content.Find.Execute(
    "En Puertollano a  31  de mayo  2022",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "En Puertollano a  26  de mayo  2022", 2)
